$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "kávéoldal" task row (row 4, columns B-F)
$ws.Range("B4").Value = "kávéoldal"
$ws.Range("C4:D4").NumberFormat = "m/d/yy h:mm"
$ws.Range("C4").Value = Get-Date -Year 2023 -Month 1 -Day 31 -Hour 17 -Minute 30 -Second 0
$ws.Range("D4").Value = Get-Date -Year 2023 -Month 1 -Day 31 -Hour 19 -Minute 0 -Second 0
$ws.Range("E4").Value = " kezdetleges"
$ws.Range("F4").Value = "elkezdtem a kávéoldalt, a háttérkép jó, a többin még dolgozom"

# Reflect the updated selection seen when the file was last saved
$ws.Range("F9").Select()

$wb.Save()
